$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "90.307.39"
$ws.Range("E2").Value = "  -0.01%  "
# Row 3
$ws.Range("D3").Value = "3.084.50"
$ws.Range("E3").Value = "  -1.59%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.92%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "624.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.34%  "
# Row 7
$ws.Range("E7").Value = "  -2.22%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.364"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.30%  "
# Row 9
$ws.Range("E9").Value = "  +0.16%  "
# Row 10
$ws.Range("D10").Value = "3.083.86"
$ws.Range("E10").Value = "  -1.44%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.724"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.85%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.196"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.18%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.95%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.49%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.54%  "
# Row 16
$ws.Range("D16").Value = "89.917.29"
$ws.Range("E16").Value = "  -0.25%  "
# Row 17
$ws.Range("D17").Value = "3.664.90"
$ws.Range("E17").Value = "  -1.26%  "
# Row 18
$ws.Range("D18").Value = "3.117.11"
$ws.Range("E18").Value = "  -0.05%  "
# Row 19
$ws.Range("E19").Value = "  +4.98%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000215"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.54%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.53%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "435.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.20%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.80%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
# Row 25
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.07%  "
# Row 26
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.31%  "
# Row 27
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.92%  "
# Row 28
$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "83.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "
# Row 29
$ws.Range("E29").Value = "  -1.89%  "
# Row 30
$ws.Range("E30").Value = "  -0.04%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.90%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.159"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "
# Row 33
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.59%  "
# Row 34
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.195"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.15%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.153"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.14%  "
# Row 36
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.88%  "
# Row 37
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.34%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "505.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.48%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.43%  "
# Row 40
$ws.Range("E40").Value = "  +0.54%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.64%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0878"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "
# Row 43
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "
# Row 44
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.407"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "
# Row 45
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +54.53%  "
# Row 46
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "151.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.21%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.685"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.24%  "
# Row 50
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.38%  "
# Row 51
$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.12%  "
